$wb = $excel.ActiveWorkbook

# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H6").Value = 56.75
$ws.Range("I6").Value = 51.88889
$ws.Range("J6").Value = 71.333336
$ws.Range("K6").Value = 155.66667
$ws.Range("L6").Value = 214.000008
$ws.Range("M6").Value = -43.66667000000001
$ws.Range("N6").Value = -438.000008
$ws.Range("H40").Value = 3481.4482
$ws.Range("I40").Value = 2430.842
$ws.Range("J40").Value = 5477.6
$ws.Range("K40").Value = 2430.842
$ws.Range("L40").Value = 5477.6
$ws.Range("M40").Value = -2255.842
$ws.Range("N40").Value = -5827.6
$ws.Range("H51").Value = 3950
$ws.Range("J51").Value = 3950
$ws.Range("L51").Value = 3950
$ws.Range("N51").Value = -4918
$ws.Range("H58").Value = 1273.4615
$ws.Range("J58").Value = 2499.1667
$ws.Range("L58").Value = 7497.500100000001
$ws.Range("N58").Value = -7797.500100000001
$ws.Range("H106").Value = 0
$ws.Range("I106").Value = 0
$ws.Range("K106").Value = 0
$ws.Range("M106").ClearContents()
$ws.Range("H113").Value = 1483.25
$ws.Range("I113").Value = 1486.3636
$ws.Range("K113").Value = 1486.3636
$ws.Range("M113").Value = 1767.6364
$ws.Range("H132").Value = 11639.228
$ws.Range("I132").Value = 11303.2
$ws.Range("J132").Value = 14999.5
$ws.Range("K132").Value = 33909.60000000001
$ws.Range("L132").Value = 44998.5
$ws.Range("M132").Value = -31379.60000000001
$ws.Range("N132").Value = -50058.5
$ws.Range("H138").Value = 1624.8
$ws.Range("J138").Value = 2199.1
$ws.Range("L138").Value = 6597.299999999999
$ws.Range("N138").Value = -16877.3

# --- Sheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H6").Value = 7500275
$ws.Range("I6").Value = 7143171.5
$ws.Range("K6").Value = 7143171.5
$ws.Range("M6").Value = -7142998.5
$ws.Range("H45").Value = 3386.95
$ws.Range("I45").Value = 2555.7856
$ws.Range("J45").Value = 5326.3335
$ws.Range("K45").Value = 2555.7856
$ws.Range("L45").Value = 5326.3335
$ws.Range("M45").Value = -2178.7856
$ws.Range("N45").Value = -6080.3335
$ws.Range("H47").Value = 67041
$ws.Range("J47").Value = 67041
$ws.Range("L47").Value = 67041
$ws.Range("M47").NumberFormat = "General"
$ws.Range("N47").Value = -68491
$ws.Range("H61").Value = 2169.8333
$ws.Range("I61").Value = 2169.8333
$ws.Range("K61").Value = 2169.8333
$ws.Range("M61").Value = -1957.8333
$ws.Range("H63").Value = 2033
$ws.Range("I63").Value = 2049.5
$ws.Range("J63").Value = 2000
$ws.Range("K63").Value = 2049.5
$ws.Range("L63").Value = 2000
$ws.Range("M63").Value = -1363.5
$ws.Range("N63").Value = -3372
$ws.Range("H66").Value = 2033
$ws.Range("I66").Value = 2049.5
$ws.Range("J66").Value = 2000
$ws.Range("K66").Value = 10247.5
$ws.Range("L66").Value = 10000
$ws.Range("M66").Value = -6815.5
$ws.Range("N66").Value = -16864
$ws.Range("H113").Value = 0
$ws.Range("J113").Value = 0
$ws.Range("L113").Value = 0
$ws.Range("N113").ClearContents()
$ws.Range("H122").Value = 1268
$ws.Range("I122").Value = 1085
$ws.Range("J122").Value = 2000
$ws.Range("K122").Value = 3255
$ws.Range("L122").Value = 6000
$ws.Range("M122").Value = -805
$ws.Range("N122").Value = -10900
$ws.Range("H136").Value = 2169.8333
$ws.Range("I136").Value = 2169.8333
$ws.Range("K136").Value = 6509.499899999999
$ws.Range("M136").Value = -3959.499899999999

# --- Sheet: BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H16").Value = 1233.3334
$ws.Range("I16").Value = 1233.3334
$ws.Range("K16").Value = 1233.3334
$ws.Range("M16").Value = -1063.3334
$ws.Range("H94").Value = 947.25
$ws.Range("I94").Value = 930
$ws.Range("J94").Value = 999
$ws.Range("K94").Value = 930
$ws.Range("L94").Value = 999
$ws.Range("M94").Value = -479
$ws.Range("N94").Value = -1901
$ws.Range("H95").Value = 16103.667
$ws.Range("J95").Value = 16103.667
$ws.Range("L95").Value = 16103.667
$ws.Range("N95").Value = -21595.667

# --- Sheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 6413.206
$ws.Range("J31").Value = 9461.789000000001
$ws.Range("L31").Value = 9461.789000000001
$ws.Range("N31").Value = -10051.789
$ws.Range("H34").Value = 6413.206
$ws.Range("J34").Value = 9461.789000000001
$ws.Range("L34").Value = 9461.789000000001
$ws.Range("N34").Value = -9865.789000000001
$ws.Range("H94").Value = 6258.5
$ws.Range("I94").Value = 2345.6667
$ws.Range("K94").Value = 2345.6667
$ws.Range("M94").Value = -1894.6667
$ws.Range("H103").Value = 16368.8
$ws.Range("I103").Value = 16368.8
$ws.Range("K103").Value = 16368.8
$ws.Range("M103").Value = -15196.8
$ws.Range("H105").Value = 1516.5
$ws.Range("J105").Value = 2148
$ws.Range("L105").Value = 2148
$ws.Range("N105").Value = -5642

# --- Sheet: CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H38").Value = 297.08
$ws.Range("I38").Value = 358.72223
$ws.Range("J38").Value = 138.57143
$ws.Range("K38").Value = 1076.16669
$ws.Range("L38").Value = 415.71429
$ws.Range("M38").Value = -729.16669
$ws.Range("N38").Value = -1109.71429

# --- Sheet: GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H2").Value = 114.066666
$ws.Range("I2").Value = 56
$ws.Range("J2").Value = 152.77777
$ws.Range("K2").Value = 56
$ws.Range("L2").Value = 152.77777
$ws.Range("M2").Value = 57
$ws.Range("N2").Value = -378.77777
$ws.Range("H10").Value = 7667
$ws.Range("I10").Value = 7667
$ws.Range("K10").Value = 7667
$ws.Range("M10").Value = -7498
$ws.Range("H14").Value = 8689.714
$ws.Range("J14").Value = 12086
$ws.Range("L14").Value = 12086
$ws.Range("N14").Value = -12422
$ws.Range("H18").Value = 8668.333000000001
$ws.Range("I18").Value = 8668.333000000001
$ws.Range("K18").Value = 8668.333000000001
$ws.Range("M18").Value = -8375.333000000001
$ws.Range("H29").Value = 17000
$ws.Range("J29").Value = 17000
$ws.Range("L29").Value = 17000
$ws.Range("N29").Value = -17580
$ws.Range("H33").Value = 9244.875
$ws.Range("I33").Value = 7000
$ws.Range("J33").Value = 9993.166999999999
$ws.Range("K33").Value = 7000
$ws.Range("L33").Value = 9993.166999999999
$ws.Range("M33").Value = -6748
$ws.Range("N33").Value = -10497.167
$ws.Range("H35").Value = 28344.666
$ws.Range("I35").Value = 9000
$ws.Range("K35").Value = 9000
$ws.Range("M35").Value = -8702
$ws.Range("H36").Value = 2814.2856
$ws.Range("I36").Value = 1166.6666
$ws.Range("J36").Value = 4050
$ws.Range("K36").Value = 1166.6666
$ws.Range("L36").Value = 4050
$ws.Range("M36").Value = -681.6666
$ws.Range("N36").Value = -5020
$ws.Range("H43").Value = 18000
$ws.Range("J43").Value = 18000
$ws.Range("L43").Value = 18000
$ws.Range("N43").Value = -18302
$ws.Range("H48").Value = 18000
$ws.Range("J48").Value = 18000
$ws.Range("L48").Value = 18000
$ws.Range("N48").Value = -18970
$ws.Range("H93").Value = 26999.75
$ws.Range("J93").Value = 25856.857
$ws.Range("L93").Value = 25856.857
$ws.Range("N93").Value = -29600.857
$ws.Range("H102").Value = 4754
$ws.Range("I102").Value = 4449.143
$ws.Range("K102").Value = 4449.143
$ws.Range("M102").Value = -2827.143

# --- Sheet: LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H16").Value = 813.4286
$ws.Range("I16").Value = 398.33334
$ws.Range("K16").Value = 398.33334
$ws.Range("M16").Value = -228.33334
$ws.Range("H20").Value = 379499.75
$ws.Range("I20").Value = 172666.67
$ws.Range("J20").Value = 999999
$ws.Range("K20").Value = 172666.67
$ws.Range("L20").Value = 999999
$ws.Range("M20").Value = -172440.67
$ws.Range("N20").Value = -1000451
$ws.Range("H40").Value = 7155.579
$ws.Range("I40").Value = 6748.25
$ws.Range("J40").Value = 7853.857
$ws.Range("K40").Value = 6748.25
$ws.Range("L40").Value = 7853.857
$ws.Range("M40").Value = -6612.25
$ws.Range("N40").Value = -8125.857
$ws.Range("H42").Value = 9995
$ws.Range("I42").Value = 9995
$ws.Range("K42").Value = 9995
$ws.Range("M42").Value = -9432
$ws.Range("H49").Value = 9995
$ws.Range("I49").Value = 9995
$ws.Range("K49").Value = 9995
$ws.Range("M49").Value = -9848
$ws.Range("H93").Value = 1802.5454
$ws.Range("I93").Value = 1800.75
$ws.Range("J93").Value = 1803.5714
$ws.Range("K93").Value = 1800.75
$ws.Range("L93").Value = 1803.5714
$ws.Range("M93").Value = -552.75
$ws.Range("N93").Value = -4299.5714
$ws.Range("H122").Value = 3709.8572
$ws.Range("I122").Value = 3709.8572
$ws.Range("K122").Value = 11129.5716
$ws.Range("M122").Value = -8679.571599999999

# --- Sheet: WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H58").Value = 20799.834
$ws.Range("I58").Value = 12000
$ws.Range("K58").Value = 12000
$ws.Range("M58").Value = -11692
$ws.Range("H95").Value = 31999.857
$ws.Range("J95").Value = 31999.857
$ws.Range("L95").Value = 31999.857
$ws.Range("N95").Value = -37491.857
$ws.Range("H122").Value = 3051.5938
$ws.Range("I122").Value = 2006.6111
$ws.Range("K122").Value = 6019.8333
$ws.Range("M122").Value = -3569.8333
$ws.Range("H135").Value = 50595.8
$ws.Range("J135").Value = 50595.8
$ws.Range("L135").Value = 50595.8
$ws.Range("N135").Value = -60735.8
